$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, shifting existing rows 102:215 down to 103:216
$ws.Rows(102).Insert()

# Populate the newly inserted row 102 with the new data record
$ws.Cells.Item(102, 1).Value = 5
$ws.Cells.Item(102, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(102, 3).Value = "Maule"
$ws.Cells.Item(102, 4).Value = 44494
$ws.Cells.Item(102, 5).Value = 7
$ws.Cells.Item(102, 6).Value = 100112032
$ws.Cells.Item(102, 7).Value = "Zapallo italiano"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 400
$ws.Cells.Item(102, 11).Value = 11000
$ws.Cells.Item(102, 12).Value = 11000
$ws.Cells.Item(102, 13).Value = 11000
$ws.Cells.Item(102, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(102, 15).Value = "Región del Maule"
$ws.Cells.Item(102, 16).Value = 183
$ws.Cells.Item(102, 17).Value = 60
$ws.Cells.Item(102, 18).Value = "Hortaliza"
